$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Methods/Children" value (column E) needs to move into the
# "Parent" column (column D), merging methods and parents into one column.
$rows = @(3, 4, 6, 7, 9, 12, 14, 15, 19, 20, 22, 23, 25, 26, 27)

foreach ($r in $rows) {
    $srcCell = $ws.Cells.Item($r, 5)   # column E
    $dstCell = $ws.Cells.Item($r, 4)   # column D
    $dstCell.Value2 = $srcCell.Value2
    $srcCell.ClearContents()
}
